# MultiLineDate.xlsx edit
# - Adds a new "Custom string property name" column (H) with three Lorem-ipsum
#   style sample values under it (H2:H4), mirroring the other sample columns.
# - Selects F10 (new "current cell") instead of the old F5.
# - Auto-sizes the newly-touched columns (A-E, G-H) to fit their content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + sample data in column H
$ws.Range("H1").Value = "Custom string property name"
$ws.Range("H2").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit."
$ws.Range("H3").Value = "Sed do eiusmod tempor incididunt ut labore et dolore."
$ws.Range("H4").Value = "Ut enim ad minim veniam, quis nostrud exercitation."

# Resize columns to fit their (new) contents - leave column F (DateTimeTest)
# alone since it already has an explicit width.
$ws.Range("A1:E4").EntireColumn.AutoFit() | Out-Null
$ws.Range("G1:H4").EntireColumn.AutoFit() | Out-Null

# Move the active selection like the author did
$ws.Range("F10").Select() | Out-Null
